# "4.0.3 model and data"
#
# The transportation input-file lists for the Boolean (constrained-to-values)
# sheet split the single "BVTQaZ" and "VTQaZ" CSV file references into six
# separate per-vehicle-type files each (LDVs, HDVs, aircraft, rail, ships,
# motorbikes). Six blank rows are also appended at the bottom of that sheet.
# View/selection state is updated to reflect where the editor was working
# when the file was saved (About sheet active/selected; Integer sheet's
# selection parked at A13; Boolean sheet scrolled down with A32 selected).

$wb = $excel.ActiveWorkbook

$wsAbout      = $wb.Worksheets.Item("About")
$wsInteger    = $wb.Worksheets.Item("Integer")
$wsBoolean    = $wb.Worksheets.Item("Boolean")
$wsSubscript  = $wb.Worksheets.Item("Subscript")

# --- Boolean sheet: split "trans/BVTQaZ/BVTQaZ.csv" (row 17) into 6 rows ---
$wsBoolean.Rows.Item(18).Insert()
$wsBoolean.Rows.Item(18).Insert()
$wsBoolean.Rows.Item(18).Insert()
$wsBoolean.Rows.Item(18).Insert()
$wsBoolean.Rows.Item(18).Insert()

$wsBoolean.Cells.Item(17, 1).Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBoolean.Cells.Item(18, 1).Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBoolean.Cells.Item(19, 1).Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBoolean.Cells.Item(20, 1).Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBoolean.Cells.Item(21, 1).Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBoolean.Cells.Item(22, 1).Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# --- Boolean sheet: split "trans/VTQaZ/VTQaZ.csv" (now row 26) into 6 rows ---
$wsBoolean.Rows.Item(27).Insert()
$wsBoolean.Rows.Item(27).Insert()
$wsBoolean.Rows.Item(27).Insert()
$wsBoolean.Rows.Item(27).Insert()
$wsBoolean.Rows.Item(27).Insert()

$wsBoolean.Cells.Item(26, 1).Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBoolean.Cells.Item(27, 1).Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBoolean.Cells.Item(28, 1).Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBoolean.Cells.Item(29, 1).Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBoolean.Cells.Item(30, 1).Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBoolean.Cells.Item(31, 1).Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# --- Boolean sheet: six new blank rows appended after the last data row (32) ---
$wsBoolean.Rows.Item(33).Insert()
$wsBoolean.Rows.Item(33).Insert()
$wsBoolean.Rows.Item(33).Insert()
$wsBoolean.Rows.Item(33).Insert()
$wsBoolean.Rows.Item(33).Insert()
$wsBoolean.Rows.Item(33).Insert()

# --- View/selection bookkeeping, in the order the sheets were last visited ---
$wsInteger.Activate()
$wsInteger.Range("A13").Select()

$wsBoolean.Activate()
$wsBoolean.Application.ActiveWindow.ScrollRow = 10
$wsBoolean.Range("A32").Select()

$wsSubscript.Activate()

$wsAbout.Activate()
